# Update Income details sheet:
#  - Insert three new income rows (Project, Stocks, Investments) above the
#    existing "Side Project" row, pushing it from row 2 down to row 5.
#  - Update the (now shifted) "Side Project" row's Amount and Date values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before the existing row 2 ("Side Project"),
# shifting it (and its formatting) down to row 5.
$ws.Rows("2:4").Insert()

# Propagate the date cell formatting (numFmtId 14) from the shifted
# "Side Project" row down into the newly inserted rows' C column.
$ws.Range("C5").Copy()
$ws.Range("C2:C4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 2: Project
$ws.Range("A2").Value = "Project"
$ws.Range("B2").Value = 2500
$ws.Range("C2").Value = 45884.22928240741

# Row 3: Stocks
$ws.Range("A3").Value = "Stocks"
$ws.Range("B3").Value = 3000
$ws.Range("C3").Value = 45884.22928240741

# Row 4: Investments
$ws.Range("A4").Value = "Investments"
$ws.Range("B4").Value = 5000
$ws.Range("C4").Value = 45884.22928240741

# Row 5: Side Project (previously row 2) - update Amount and Date
$ws.Range("B5").Value = 5000
$ws.Range("C5").Value = 45882.22928240741

# Expand the "numbers stored as text" ignored-error marker to cover the
# full, now-larger data range.
$ws.Range("A1:C5").Errors.Item(9).Value = $true
